$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: new trailing columns (PriceChange, UpDown)
$ws.Range("X31").Value = -0.21000000000000085
$ws.Range("Y31").Value = "Down"

# Row 32: new row of repeater data (4th bio stock)
$ws.Range("A32").Value = 42651.425138888888
$ws.Range("B32").Value = 13
$ws.Range("C32").Value = "Buy"
$ws.Range("D32").Value = 58
$ws.Range("E32").Value = 1323
$ws.Range("F32").Value = 76
$ws.Range("G32").Value = 70
$ws.Range("H32").Value = 29
$ws.Range("I32").Value = 92
$ws.Range("J32").Value = 7
$ws.Range("K32").Value = 13326
$ws.Range("L32").Value = 12
$ws.Range("M32").Value = 5
$ws.Range("N32").Value = 13
$ws.Range("O32").Value = 1
$ws.Range("P32").Value = "Named"
$ws.Range("Q32").Value = 47.963765586266284
$ws.Range("R32").Value = 0.49
$ws.Range("S32").Value = 0.0521
$ws.Range("T32").Value = -0.0214
$ws.Range("U32").Value = 2.2799999999999998
$ws.Range("V32").Value = "N/A"
$ws.Range("W32").Value = 0

# Date column (A) uses date format style; Percent columns (S,T) use percent style
$ws.Range("A32").NumberFormat = "m/d/yy h:mm"
$ws.Range("S32").NumberFormat = "0.00%"
$ws.Range("T32").NumberFormat = "0.00%"
